$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": Right marks 5 -> 4, Wrong marks -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total 100 -> 80, Wrong total -1 -> -2, and the
# "obtained / max" summary text updated accordingly
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "78 / 112"
